# Append a new lancers.jp listing to the "ランサーズ" sheet (2026-01-25 12:51 JST run),
# pushing the previous rows 7-9 down to 8-10, and refresh the collection timestamp
# (column A) on every data row to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-25 12:51:08"

# Insert a new row at position 7; existing rows 7-9 (and their formatting) shift down to 8-10.
$ws.Rows(7).Insert()

# Populate the newly inserted row 7 with the new listing's data.
$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "【3名限定】10万円でアプリ開発をして欲しい方探しています!"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5478575"
$ws.Range("G7").Value = 85
$ws.Range("H7").Value = "◆開発 ◇アプリ"

# Refresh the "取得日時" (retrieved at) timestamp for every data row (2-10) to the new run time.
for ($i = 2; $i -le 10; $i++) {
    $ws.Range("A" + $i).Value = $newTimestamp
}

# Row.Insert() does not carry hyperlink targets along with the shifted cells, so rebuild the
# hyperlinks for column F (rows 2-10) so each one again points at the URL shown in its own cell.
$ws.Hyperlinks.Delete()
for ($i = 2; $i -le 10; $i++) {
    $cell = $ws.Range("F" + $i)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}
